$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7..87 down to 8..88
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new record's data
$ws.Cells.Item(7,1).Value = 10
$ws.Cells.Item(7,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(7,3).Value = "La Araucanía"
$ws.Cells.Item(7,4).Value = 44515
$ws.Cells.Item(7,5).Value = 9
$ws.Cells.Item(7,6).Value = 100112031
$ws.Cells.Item(7,7).Value = "Poroto verde"
$ws.Cells.Item(7,8).Value = "Sin especificar"
$ws.Cells.Item(7,9).Value = "Primera"
$ws.Cells.Item(7,10).Value = 45
$ws.Cells.Item(7,11).Value = 45000
$ws.Cells.Item(7,12).Value = 47000
$ws.Cells.Item(7,13).Value = 45889
$ws.Cells.Item(7,14).Value = "$/malla 25 kilos"
$ws.Cells.Item(7,15).Value = "Perú"
$ws.Cells.Item(7,16).Value = 1836
$ws.Cells.Item(7,17).Value = 25
$ws.Cells.Item(7,18).Value = "Hortaliza"
